$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the crypto snapshot: coinranking.com scrape, updated by the
# scheduled GitHub Actions job. Columns D (Price) and E (Volume/1h) are
# always stored as plain text in this sheet (values like "35.077.90" or
# "  +1.73%  " are display strings, not numbers/percentages). For any new
# value that also happens to parse as a plain number (e.g. "237.27"),
# force the cell to Text format first so it is written back as a literal
# string instead of being coerced into a numeric cell -- then restore the
# default "Normal" style so no stray formatting is left behind.

$ws.Range('D2').Value2 = '35.077.90'
$ws.Range('E2').Value2 = '  +1.73%  '
$ws.Range('D3').Value2 = '1.856.66'
$ws.Range('E3').Value2 = '  +3.17%  '
$ws.Range('E4').Value2 = '  +0.23%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value2 = '237.27'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value2 = '  +3.59%  '
$ws.Range('E6').Value2 = '  +1.75%  '
$ws.Range('E7').Value2 = '  +0.14%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value2 = '42.44'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value2 = '  +8.17%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value2 = '0.329'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value2 = '  +3.06%  '
$ws.Range('E11').Value2 = '  +0.23%  '
$ws.Range('D12').Value2 = '2.126.11'
$ws.Range('E12').Value2 = '  +3.24%  '
$ws.Range('D13').Value2 = '1.862.13'
$ws.Range('E13').Value2 = '  +3.78%  '
$ws.Range('E14').Value2 = '  +3.14%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value2 = '0.678'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value2 = '  +3.27%  '
$ws.Range('E16').Value2 = '  +3.43%  '
$ws.Range('D17').Value2 = '35.056.54'
$ws.Range('E17').Value2 = '  +2.13%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value2 = '70.31'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value2 = '  +1.98%  '
$ws.Range('D19').Value2 = '0.0₃0796'
$ws.Range('E19').Value2 = '  +2.39%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value2 = '240.73'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value2 = '  +0.68%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value2 = '12.14'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value2 = '  +3.28%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value2 = '4.74'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value2 = '  +1.51%  '
$ws.Range('E23').Value2 = '  +0.08%  '
$ws.Range('E24').Value2 = '  +1.91%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value2 = '171.19'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value2 = '  -1.09%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value2 = '1.90'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value2 = '  +27.79%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value2 = '7.93'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value2 = '  +3.19%  '
$ws.Range('E28').Value2 = '  +3.13%  '
$ws.Range('E29').Value2 = '  +2.10%  '
$ws.Range('E30').Value2 = '  +0.21%  '
$ws.Range('E31').Value2 = '  +3.22%  '
$ws.Range('E32').Value2 = '  +0.55%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value2 = '4.03'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value2 = '  +3.69%  '
$ws.Range('E34').Value2 = '  +13.20%  '
$ws.Range('E35').Value2 = '  +22.84%  '
$ws.Range('B36').Value2 = 'TrustWalletToken'
$ws.Range('C36').Value2 = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value2 = '1.29'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value2 = '  +5.62%  '
$ws.Range('B37').Value2 = 'ImmutableX'
$ws.Range('C37').Value2 = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value2 = '0.782'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value2 = '  +12.78%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value2 = '91.85'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value2 = '  +1.60%  '
$ws.Range('E40').Value2 = '  +7.22%  '
$ws.Range('D41').Value2 = '1.352.73'
$ws.Range('E41').Value2 = '  +2.22%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value2 = '14.86'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value2 = '  +5.44%  '
$ws.Range('E43').Value2 = '  +5.96%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value2 = '12.77'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value2 = '  +54.47%  '
$ws.Range('E45').Value2 = '  +0.83%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value2 = '2.74'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value2 = '  +1.59%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value2 = '0.0543'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value2 = '  +6.19%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value2 = '6.44'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value2 = '  +4.92%  '
$ws.Range('D49').Value2 = '2.039.96'
$ws.Range('E49').Value2 = '  +2.92%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value2 = '0.0681'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value2 = '  +3.30%  '
$ws.Range('E51').Value2 = '  +18.25%  '
